$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 5 with the new test-mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Ik heb nog geen geld terug."
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #3: Ik heb nog geen geld terug."
$logs.Range("D5").Value = "Retour / Terugbetaling"
$logs.Range("E5").Value = "Dank voor je bericht. We hebben je eerdere e-mail ontvangen en doorgestuurd naar retour@bedrijf.nl."
$logs.Range("F5").Value = "2025-08-04 20:03:35"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Ja"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# Extend the conditional-formatting blocks from row 4 to row 5
$logs.Range("D2:D4").FormatConditions().Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions().Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions().Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions().Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions().Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# --- Dashboard sheet: append row 3 with the new category count ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Retour / Terugbetaling"
$dash.Range("B3").Value = 1

# --- Chart: extend the category/value series references to include row 3 ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart()
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
